$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Insert a new "task_details" row right before the last task block
#    (old row 32) so it becomes an extra Mission/Responsibility row
#    appended to the previous task's details table.
# ------------------------------------------------------------------
$ws.Rows("32").Insert()

$ws.Range("F32").Value = "bvcbvcbvcbvcbvcsgfgsd"
$ws.Range("G32").Value = "cbvcbcvbvc"
$ws.Range("H32").Value = "bvcbvcbcv"
$ws.Range("I32").Value = "Saab@gfgfd.com"
$ws.Range("J32").Value = "Stuck"

# ------------------------------------------------------------------
# 2) Style all the task_details header rows (the "Mission / Responsibility
#    / email / status" rows) with bold font + red fill.
# ------------------------------------------------------------------
$headerRows = @(3, 14, 17, 22, 28, 34)
foreach ($r in $headerRows) {
    $rng = $ws.Range("F" + $r + ":J" + $r)
    $rng.Font.Bold = $true
    $rng.Interior.Color = 8421631
}

# ------------------------------------------------------------------
# 3) Append two more task_details rows to the new task's table
#    (rows 39 and 40, right after the existing row 38).
# ------------------------------------------------------------------
$ws.Range("F39").Value = "cxzcxzcxz"
$ws.Range("G39").Value = "cxzcxzcxz"
$ws.Range("H39").Value = "DUDU"
$ws.Range("I39").Value = "benharushtomer@gmail.com"
$ws.Range("J39").Value = "Closed"

$ws.Range("F40").Value = "vcxvcxvcx"
$ws.Range("G40").Value = "vcxvcxvcx"
$ws.Range("I40").Value = "vcxvcx@gmail.com"
$ws.Range("J40").Value = "Closed"

Write-Output "done"
